$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (previously 005411104 / PAULO / 233000) becomes the SABRINA row,
# replacing what used to be three rows (PAULO, LOHRAN, MARIANA).
$ws.Cells.Item(2, 1).NumberFormat = "@"
$ws.Cells.Item(2, 1).Value = "004301699"
$ws.Cells.Item(2, 1).ClearFormats()
$ws.Cells.Item(2, 2).Value = "SABRINA"
$ws.Cells.Item(2, 3).Value = 25433.49

# Remove the old LOHRAN (row 3) and MARIANA (row 4) rows entirely -
# THIAGO (was row 5) shifts up to row 3.
$ws.Rows.Item(4).Delete()
$ws.Rows.Item(3).Delete()

# Remove the FERNANDO and CAROLINE rows that followed THIAGO.
$ws.Rows.Item(5).Delete()
$ws.Rows.Item(4).Delete()
